# Daily attendance processing - 2025-11-13 18:56:18
# Reorders the "Recorded By" (column G) author lists so that the
# automated "System" entry is listed after the human/backup recorder
# instead of before it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Known exact replacements observed for the "Recorded By" column (G).
# Only rows whose current text matches one of these keys are touched;
# everything else (single values, admin@admin.com combos, values that
# already list System last, etc.) is left exactly as-is.
$replacements = @{
    'System, backup@backdoor.com'           = 'backup@backdoor.com, System'
    'System, dnasr281@gmail.com'            = 'dnasr281@gmail.com, System'
    'system, System, backup@backdoor.com'   = 'backup@backdoor.com, System, system'
}

# Column G holds "Recorded By"; data rows run from row 2 through the
# last used row on the sheet.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = $ws.UsedRange.Rows.Count }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $current = $cell.Value2
    if ($null -ne $current -and $replacements.ContainsKey($current)) {
        $cell.Value2 = $replacements[$current]
    }
}
